# Apply the changes described by the commit "avoid double load of i18n.xsl":
#
# 1) The lone paragraph in the body was carrying the "Note" paragraph
#    style for no good reason (it only holds the _GoBack bookmark) -
#    clear it back to the document default (Normal) so the <w:pPr>
#    wrapper with <w:pStyle w:val="Note"/> disappears entirely.
#
# 2) Fix the typo in the custom style name/id "MarginNoteRIght" ->
#    "MarginNoteRight" (note the lower-case "i") everywhere it is used.

$d = $word.ActiveDocument

# --- 1) Drop the stray "Note" paragraph style from the first paragraph ---
$p = $d.Paragraphs(1)
if ($p.Range.Style.NameLocal -eq "Note") {
    $p.Range.Style = $d.Styles("Normal")
}

# --- 2) Rename the misspelled "MarginNoteRIght" style to "MarginNoteRight" ---
$style = $d.Styles("MarginNoteRIght")
$style.NameLocal = "MarginNoteRight"
